# Add a new affiliation row (row 7) to the worksheet:
#   Association for Information Science and Technology | -
#   | Silver Spring, MD | https://www.asist.org/
# and carry the running index formula (=A6+1) down into A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data.
$ws.Range("A7").Formula = "=A6+1"
$ws.Range("B7").Value = "Association for Information Science and Technology"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "Silver Spring, MD"
$ws.Range("E7").Value = "https://www.asist.org/"

# Matches the author's final cursor position in the saved file.
[void]$ws.Range("E7").Select()
